$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: agravo1.pdf -> agravo3.pdf, CNJ value changes
$ws.Range("A2").Value = "agravo3.pdf"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "0001547-36.2008.4.01.3814"

# Remove rows 3 through 5 (old agravo2, agravo3, agravo4 entries) entirely
$ws.Range("A3:C5").EntireRow.Delete() | Out-Null
